$d = $word.ActiveDocument

# 1) Drop the old "_GoBack" bookmark that currently sits inside
#    "titel = varchar(100)" (between the "10" and "0)" runs). Word only
#    ever keeps a single "_GoBack" bookmark, and its new home is the
#    "e-mailadres" paragraph edited below, so just remove it here first
#    (cleanly, via the native Bookmarks collection, before any XML
#    surgery touches range offsets).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) "e-mailadres = varchar(255)" paragraph: append a new, separate
#    " Unique" run (en-GB language) and re-create the "_GoBack" bookmark
#    right after it, at the end of the paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*e-mailadres = varchar(255)*") {
        $full = $p.Range
        $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='169FE0FB' w14:textId='77777777' w:rsidR='00B34C07' w:rsidRPr='004E78D7' w:rsidRDefault='00B34C07' w:rsidP='00B34C07'>" +
               "<w:pPr><w:pStyle w:val='NoSpacing'/><w:rPr><w:lang w:val='en-GB'/></w:rPr></w:pPr>" +
               "<w:r w:rsidRPr='004E78D7'><w:rPr><w:lang w:val='en-GB'/></w:rPr><w:t>e-mailadres = varchar(255)</w:t></w:r>" +
               "<w:r><w:rPr><w:lang w:val='en-GB'/></w:rPr><w:t xml:space='preserve'> Unique</w:t></w:r>" +
               "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
               "</w:p>"
        $full.InsertXML($xml)
        break
    }
}
